$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value2 = '7EtlCyYSj138nJcU41wo'
$ws.Range("D2").Value2 = 'iT7LGm2F06jcT3zHNzit'
$ws.Range("E2").Value2 = 'IH8LR8gnCyYoVq5dd2fg'
$ws.Range("F2").Value2 = '4q5Epfo0jW73PsloHgR5'
$ws.Range("C3").Value2 = 'jXvEwk5NYvCjfIMRZO3T'
$ws.Range("D3").Value2 = '9eiKUqXrfWnWMs41gtZX'
$ws.Range("E3").Value2 = 'WbPYnchpdBiF6T6sDCQP'
$ws.Range("F3").Value2 = '8Pb6SqV3IjyTX5xseUkL'
$ws.Range("C4").Value2 = 'i0ZAiLD4lr77kLbtCOXl'
$ws.Range("D4").Value2 = 'Ji6mIQ8vwj6Ht6SZHJtE'
$ws.Range("E4").Value2 = '7m4MWRpzlX6rQCzFrWtY'
$ws.Range("F4").Value2 = 'Vaz2PGhGgS7VB73icss2'
$ws.Range("C5").Value2 = 'FSBarAK6DvDh9g4qvaxv'
$ws.Range("D5").Value2 = 'drmFfKInneT6yPpElwyi'
$ws.Range("E5").Value2 = 'ETZI9rLfsyc6Z1goYFYe'
$ws.Range("F5").Value2 = 'yewUNK29HeBPRR1WWTG8'
$ws.Range("C6").Value2 = 'PAa99UmkpnmUfw5toqwY'
$ws.Range("D6").Value2 = 'djqFV46PW8rtZQRCLMDS'
$ws.Range("C8").Value2 = 'WMx3zZ6AVRNrHX2UgAAu'
$ws.Range("D8").Value2 = 'rY6WjfffjgkoT7w2ilwZ'
$ws.Range("E8").Value2 = 'Boa7tH8KtCSbecskjToJ'
$ws.Range("F8").Value2 = 'cDpl3AD7szgjb8KCg1QS'
$ws.Range("C9").Value2 = '/Kw45gpMP8CXoZRiFyGQo'
$ws.Range("D9").Value2 = 'wmcVggFs432aKn0fsKZ2'
$ws.Range("E9").Value2 = '3uocvVlBY6Uoog97bp2k'
$ws.Range("F9").Value2 = 'JfKq84Xmkuq66dsJb7Uy'
$ws.Range("C10").Value2 = 'DL4ULiomK4Ndwn5qRuzQ'
$ws.Range("D10").Value2 = 'VcTplGlDfq62BdCTmH99'
$ws.Range("E10").Value2 = 'yKcOVDRP33qLekAiiGQz'
$ws.Range("F10").Value2 = 'A3huw8jZ9CQ9EZ582CHl'
$ws.Range("C11").Value2 = 'YjZzvIKaLrQEqMSAQTXa'
$ws.Range("D11").Value2 = 'P7sNr7j6GOU1zasc5JE4'
$ws.Range("E11").Value2 = 'oMU92DQUcdbQ3LxxOKlj'
$ws.Range("F11").Value2 = 'AeX8c1Ar0ly1vBaMDXja'
$ws.Range("C12").Value2 = 'tfYZikokPbPzQidwLt1b'
$ws.Range("D12").Value2 = 'R79ypq2z5STFoqlsKsAD'
$ws.Range("C14").Value2 = 'mXv1pOSkB3NCTVVlHraK'
$ws.Range("D14").Value2 = 'nxls50wwacQV9Aq54wfH'
$ws.Range("E14").Value2 = 'CWnDqASaBBSiCu5Oq74Q'
$ws.Range("F14").Value2 = 'gj7EbiobLR83wvEvuumT'
$ws.Range("C15").Value2 = 'HVsCszS8NnqVbRqX6KjA'
$ws.Range("D15").Value2 = '8Z4nTC7O5twUtkJTfDzF'
$ws.Range("E15").Value2 = 'GPzgoJO65cyJFlOk6w5m'
$ws.Range("F15").Value2 = 'AuPOQvm13KhRRy3NJIqN'
$ws.Range("C16").Value2 = 'lbs6iRXDzTQUFxyYtXkA'
$ws.Range("D16").Value2 = 'vyRoCdsoOFhC7QeAVpWH'
$ws.Range("E16").Value2 = 'KS4RUJ3pHE65Wn9vkTvD'
$ws.Range("F16").Value2 = 'GJczcgHgbBQ5r9tJSCA9'
$ws.Range("C17").Value2 = 'vVdViHEjcvwidczMIv8s'
$ws.Range("D17").Value2 = 'GFkETYHUSpGBEL4r7xmw'
$ws.Range("E17").Value2 = '26lvgplRN5MqnNCw49lP'
$ws.Range("F17").Value2 = '0Z9LY8oZS9z1wS2vabhN'
$ws.Range("C18").Value2 = 'sviTU4fJcWW0gqSALlgn'
$ws.Range("D18").Value2 = 'XpPBSE4pCOgqEIfFpKew'
$ws.Range("C20").Value2 = 'VKIOiSHUqyVld2Yq1V14'
$ws.Range("D20").Value2 = '1XHw3PnaJFnLdwG7Yt8M'
$ws.Range("E20").Value2 = 'VRAX8v7tJ0ihLxHiH18C'
$ws.Range("F20").Value2 = 'pItXz2bFTKYR6jEzfYMg'
$ws.Range("C21").Value2 = 'Te53vp3IGZaQArIu9STR'
$ws.Range("D21").Value2 = 'V9elgbW5PbgZk9BNnbCK'
$ws.Range("E21").Value2 = '6qIBB5qLorpuGplacBQF'
$ws.Range("F21").Value2 = 'fzuowJdRUTTC5EtBx0gD'
$ws.Range("C22").Value2 = 'IpeexauMq2OIGifUYNVu'
$ws.Range("D22").Value2 = 'o42O11lXDsGXfKTgi2kE'
$ws.Range("E22").Value2 = 'qoIk4PxKAUTiw2Ky0AbK'
$ws.Range("F22").Value2 = 'rzlQ18K9OIvsxu9ZiFVR'
$ws.Range("C23").Value2 = 'Bvm6NNe7TQTRM1tfmPQ0'
$ws.Range("D23").Value2 = '5SVJ005HYbSnqjt7ffwb'
$ws.Range("E23").Value2 = 'nTvOVBLI6fz4lv95uh0y'
$ws.Range("F23").Value2 = 'fPtE2G6IqBEQbNlXBMn4'
$ws.Range("E24").Value2 = 'f2IBTdUuHj7lczMfaGB9'
$ws.Range("F24").Value2 = 'wyMhYXxkeLyyiq0qUpM9'

$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("F20").Select()
